$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.186.74'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.10%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.581.23'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.39%  '

$ws.Range('E4').Value = '  -0.27%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.11%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.497'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.14%  '

$ws.Range('E7').Value = '  -0.26%  '

$ws.Range('E8').Value = '  -1.66%  '

$ws.Range('E9').Value = '  -0.85%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.53'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.09%  '

$ws.Range('E11').Value = '  -0.30%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.803.50'
$ws.Range('D12').Style = 'Normal'

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.592.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.34%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.05'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.05%  '

$ws.Range('E15').Value = '  -1.51%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.94%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.198.45'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').Value = '0.0₃0734'
$ws.Range('E18').Value = '  -0.80%  '

$ws.Range('E19').Value = '  +1.29%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '207.70'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.36%  '

$ws.Range('E22').Value = '  -0.56%  '

$ws.Range('E23').Value = '  -2.87%  '

$ws.Range('E24').Value = '  -1.40%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.57%  '

$ws.Range('E26').Value = '  -0.23%  '

$ws.Range('E27').Value = '  -1.56%  '

$ws.Range('E28').Value = '  -1.53%  '

$ws.Range('E29').Value = '  -1.22%  '

$ws.Range('E30').Value = '  -1.57%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.85%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.88%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.95'
$ws.Range('D33').Style = 'Normal'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.275.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.51%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.46'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.31%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.611'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.67%  '

$ws.Range('E37').Value = '  -1.05%  '

$ws.Range('E38').Value = '  -1.93%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.817'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.93%  '

$ws.Range('E40').Value = '  -12.14%  '

$ws.Range('E41').Value = '  +2.41%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.79%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.764'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.94%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.27'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.10%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.717.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.22%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.80%  '

$ws.Range('E47').Value = '  -0.38%  '

$ws.Range('E48').Value = '  -1.28%  '

$ws.Range('E49').Value = '  -2.21%  '

$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.61%  '

$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.01%  '
